$d = $word.ActiveDocument

# The "date of issue" cell in the diploma-supplement header currently reads
# "     30.01.#CurrentYear" (a leading run of five spaces followed by three
# runs spelling out a fixed day/month and a #CurrentYear placeholder, each
# carrying an explicit uk-UA run language). Replace the day/month/placeholder
# portion with a single #SupplDate placeholder run, and drop the redundant
# uk-UA language override (both on the paragraph mark and on the surviving
# run) since it only duplicates the document's default language.
$rng = $d.Content
$found = $rng.Find.Execute("30.01.#CurrentYear", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $frag = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="001C6B9D" w:rsidRPr="00703BAF" w:rsidRDefault="002C64B2" w:rsidP="005B45BD"><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:line="275" w:lineRule="auto"/><w:ind w:right="57"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:spacing w:val="-1"/><w:sz w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="005B45BD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:color w:val="00B050"/><w:spacing w:val="-1"/><w:sz w:val="18"/></w:rPr><w:t xml:space="preserve">     </w:t></w:r><w:r w:rsidR="00EF7D57" w:rsidRPr="00703BAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:spacing w:val="-1"/><w:sz w:val="18"/></w:rPr><w:t>#SupplDate</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($frag)
}
